$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-8) get their variable columns (D, L, M, N, O, P, Q, R, S, T)
# rotated between rows, per the commit's weekly price-data refresh.
# Capture all "before" values first so the in-place rewrite doesn't clobber
# a row before it has been read.

$cols = @("D","L","M","N","O","P","Q","R","S","T")

$snapshot = @{}
for ($r = 2; $r -le 8; $r++) {
    $row = @{}
    foreach ($c in $cols) {
        $row[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $row
}

# after-row -> source before-row
$mapping = @{
    2 = 3
    3 = 8
    4 = 6
    5 = 2
    6 = 5
    7 = 7
    8 = 4
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $src[$c]
    }
}
